$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the quotes around the JSON object keys in the template strings
# that live in row 1 (these feed every row's J:O / R formulas via $ refs).
$ws.Range("K1").Value = ", name: "
$ws.Range("J1").Value = "{ id: "
$ws.Range("L1").Value = ", type: "
$ws.Range("M1").Value = ", value: "
$ws.Range("N1").Value = ", time: "
$ws.Range("O1").Value = ", season: "

# Move the active selection to O3, matching the author's resulting cursor
# position after building out the (future) bug name lookup helper.
$ws.Range("O3").Select()
